$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray "23-00001" row (row 16) -- not part of the actual data set
$ws.Rows.Item(16).Delete()

# Fill in the full location data set (latitude, longitude, elevation, depth)
# for every Sample ID row already present in column A
$ws.Range("B2").Value = 40.01601977399668
$ws.Range("C2").Value = -105.25618106909128
$ws.Range("D2").Value = 1753.2133302031987
$ws.Range("E2").Value = 147.99480199174525
$ws.Range("B3").Value = 40.01354292092485
$ws.Range("C3").Value = -105.25470464593657
$ws.Range("D3").Value = 1715.7635974635975
$ws.Range("E3").Value = 180.05712771131698
$ws.Range("B4").Value = 40.02183797470178
$ws.Range("C4").Value = -105.25399594946768
$ws.Range("D4").Value = 1690.9275499614805
$ws.Range("E4").Value = 192.6146364835849
$ws.Range("B5").Value = 40.01614856608919
$ws.Range("C5").Value = -105.25239552113064
$ws.Range("D5").Value = 1762.4123474751664
$ws.Range("E5").Value = 142.45407533726956
$ws.Range("B6").Value = 40.020431238355684
$ws.Range("C6").Value = -105.25589224835406
$ws.Range("D6").Value = 1742.6391691332663
$ws.Range("E6").Value = 146.01521455537815
$ws.Range("B7").Value = 40.01807365537656
$ws.Range("C7").Value = -105.25943306632337
$ws.Range("D7").Value = 1719.883068077513
$ws.Range("E7").Value = 178.94100658385457
$ws.Range("B8").Value = 40.01801933938867
$ws.Range("C8").Value = -105.2543351010837
$ws.Range("D8").Value = 1725.5239945842616
$ws.Range("E8").Value = 161.281284476273
$ws.Range("B9").Value = 40.0142556656749
$ws.Range("C9").Value = -105.25614613984045
$ws.Range("D9").Value = 1695.3701361203985
$ws.Range("E9").Value = 104.27472326119644
$ws.Range("B10").Value = 40.0143568079631
$ws.Range("C10").Value = -105.25624078516704
$ws.Range("D10").Value = 1724.0625545729163
$ws.Range("E10").Value = 146.49614231927407
$ws.Range("B11").Value = 40.01411814721032
$ws.Range("C11").Value = -105.25169444864349
$ws.Range("D11").Value = 1758.57276864775
$ws.Range("E11").Value = 153.959717426792
$ws.Range("B12").Value = 40.01676321192857
$ws.Range("C12").Value = -105.25636888315756
$ws.Range("D12").Value = 1755.592181515533
$ws.Range("E12").Value = 108.5249143416485
$ws.Range("B13").Value = 40.0154276954004
$ws.Range("C13").Value = -105.25473056288209
$ws.Range("D13").Value = 1748.5952835347541
$ws.Range("E13").Value = 151.11088434685578
$ws.Range("B14").Value = 40.0161605061043
$ws.Range("C14").Value = -105.259966614349
$ws.Range("D14").Value = 1695.6963712921606
$ws.Range("E14").Value = 189.30169152652914
$ws.Range("B15").Value = 40.016098973807054
$ws.Range("C15").Value = -105.25259209321005
$ws.Range("D15").Value = 1762.5258006718877
$ws.Range("E15").Value = 156.27872017820394

# Leave the selection where the author left it when saving
$ws.Range("C3").Select()
